# Oktober-2015.xlsx: add a plain "Sheet1" copy of the daily-data table
# (A9:K40 on "Data Harian - Table") as a new worksheet, A1:K32, and tidy
# up the source sheet (drop the stray blank A1 cell, remove the logo
# picture, reset the selection).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# The first row used to carry a totally blank A1 cell (no value, no
# style) - clear it so it stops being serialised at all.
$ws1.Range("A1").ClearContents()

# Remove the BMKG logo picture that was anchored on the source sheet.
for ($i = $ws1.Shapes.Count; $i -ge 1; $i--) {
    $ws1.Shapes.Item($i).Delete()
}

# New worksheet, positioned right after the existing one.
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "Sheet1"

# Copy the header + 31 daily rows (A9:K40) over to A1:K32 on the new
# sheet, carrying values, number formats, borders and alignment.
$src = $ws1.Range("A9:K40")
$dst = $newSheet.Range("A1")
$src.Copy($dst)

# Match the taller, auto-fit wrapped-text row height used for the data
# rows once they live on their own sheet.
$newSheet.Range("A2:K32").RowHeight = 28.8

# Selection bookkeeping matching the edited workbook: the source sheet
# scrolls down to the table and selects it, the new sheet is the active
# tab with the whole table selected.
$ws1.Activate()
$ws1.Range("A9:K40").Select()

$newSheet.Activate()
$newSheet.Range("A1:K32").Select()

Write-Host "done"
